$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; existing rows 12-81 shift down to 13-82.
$ws.Rows(12).Insert()

# Populate the newly inserted row 12 with its data.
$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = 44670
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 100112021
$ws.Range("G12").Value = "Ají"
$ws.Range("H12").Value = "Cristal"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 15500
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15750
$ws.Range("N12").Value = "`$/saco 25 kilos"
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 630
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
